$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "27.821.49"
Set-TextValue "E2" "  +6.12%  "

# Row 3
Set-TextValue "D3" "1.734.41"
Set-TextValue "E3" "  +4.49%  "

# Row 4
Set-TextValue "E4" "  -0.21%  "

# Row 5
Set-TextValue "D5" "226.88"
Set-TextValue "E5" "  +3.31%  "

# Row 6
Set-TextValue "D6" "0.5418"
Set-TextValue "E6" "  +3.09%  "

# Row 7
Set-TextValue "D7" "1.003"
Set-TextValue "E7" "  -0.19%  "

# Row 8
Set-TextValue "D8" "0.2727"
Set-TextValue "E8" "  +1.33%  "

# Row 9
Set-TextValue "D9" "0.06687"
Set-TextValue "E9" "  +4.71%  "

# Row 10
Set-TextValue "D10" "21.99"
Set-TextValue "E10" "  +6.52%  "

# Row 11
Set-TextValue "D11" "0.07746"
Set-TextValue "E11" "  +0.62%  "

# Row 12
Set-TextValue "D12" "4.671"
Set-TextValue "E12" "  +0.82%  "

# Row 13
Set-TextValue "B13" "WrappedliquidstakedEther2.0"
Set-TextValue "C13" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D13" "1.972.64"
Set-TextValue "E13" "  +4.46%  "

# Row 14
Set-TextValue "B14" "WrappedEther"
Set-TextValue "C14" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D14" "1.708.62"
Set-TextValue "E14" "  +2.57%  "

# Row 15
Set-TextValue "D15" "0.5931"
Set-TextValue "E15" "  +4.96%  "

# Row 16
Set-TextValue "D16" "0.0₅8384"
Set-TextValue "E16" "  +1.45%  "

# Row 17
Set-TextValue "D17" "68.78"
Set-TextValue "E17" "  +4.73%  "

# Row 18
Set-TextValue "D18" "27.795.37"
Set-TextValue "E18" "  +6.06%  "

# Row 19
Set-TextValue "D19" "227.49"
Set-TextValue "E19" "  +18.52%  "

# Row 20
Set-TextValue "B20" "Dai"
Set-TextValue "C20" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D20" "1.003"
Set-TextValue "E20" "  -0.16%  "

# Row 21
Set-TextValue "B21" "Uniswap"
Set-TextValue "C21" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D21" "4.799"
Set-TextValue "E21" "  +2.27%  "

# Row 22
Set-TextValue "E22" "  +3.58%  "

# Row 23
Set-TextValue "D23" "6.205"
Set-TextValue "E23" "  +3.39%  "

# Row 24
Set-TextValue "D24" "1.004"
Set-TextValue "E24" "  -0.19%  "

# Row 25
Set-TextValue "D25" "148.23"
Set-TextValue "E25" "  +1.40%  "

# Row 26
Set-TextValue "D26" "1.728"
Set-TextValue "E26" "  +13.38%  "

# Row 27
Set-TextValue "D27" "0.1244"
Set-TextValue "E27" "  +3.64%  "

# Row 28
Set-TextValue "D28" "7.488"
Set-TextValue "E28" "  +2.45%  "

# Row 29
Set-TextValue "D29" "17.02"
Set-TextValue "E29" "  +5.96%  "

# Row 30
Set-TextValue "D30" "0.05629"
Set-TextValue "E30" "  +0.26%  "

# Row 31
Set-TextValue "D31" "1.310"
Set-TextValue "E31" "  +2.26%  "

# Row 32
Set-TextValue "D32" "3.653"
Set-TextValue "E32" "  +4.50%  "

# Row 33
Set-TextValue "D33" "3.498"
Set-TextValue "E33" "  +2.47%  "

# Row 34
Set-TextValue "D34" "1.678"
Set-TextValue "E34" "  +6.35%  "

# Row 35
Set-TextValue "D35" "0.9741"
Set-TextValue "E35" "  +2.10%  "

# Row 36
Set-TextValue "D36" "2.850"
Set-TextValue "E36" "  +2.26%  "

# Row 37
Set-TextValue "D37" "2.436"
Set-TextValue "E37" "  +1.22%  "

# Row 38
Set-TextValue "D38" "0.5984"
Set-TextValue "E38" "  +4.09%  "

# Row 39
Set-TextValue "D39" "0.01671"
Set-TextValue "E39" "  +4.36%  "

# Row 40
Set-TextValue "D40" "5.924"
Set-TextValue "E40" "  -0.51%  "

# Row 41
Set-TextValue "D41" "0.8656"
Set-TextValue "E41" "  +3.32%  "

# Row 42
Set-TextValue "D42" "1.056.66"
Set-TextValue "E42" "  +2.58%  "

# Row 43
Set-TextValue "D43" "1.003"
Set-TextValue "E43" "  -0.14%  "

# Row 44
Set-TextValue "D44" "101.76"
Set-TextValue "E44" "  +0.52%  "

# Row 45
Set-TextValue "D45" "1.877.35"
Set-TextValue "E45" "  +4.35%  "

# Row 46
Set-TextValue "D46" "0.0₈115"
Set-TextValue "E46" "  +7.93%  "

# Row 47
Set-TextValue "D47" "59.78"
Set-TextValue "E47" "  +2.26%  "

# Row 48
Set-TextValue "D48" "8.281"
Set-TextValue "E48" "  +2.90%  "

# Row 49
Set-TextValue "B49" "XinFinNetwork"
Set-TextValue "C49" "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
Set-TextValue "D49" "0.06698"
Set-TextValue "E49" "  +23.45%  "

# Row 50
Set-TextValue "B50" "Mantle"
Set-TextValue "C50" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D50" "0.4431"
Set-TextValue "E50" "  +1.92%  "

# Row 51
Set-TextValue "D51" "1.001"
Set-TextValue "E51" "  +0.39%  "
